# Update user_data.xlsx:
#  - B3 "30" (text) -> 30 (number)
#  - Append row 4: A4="hi", B4="5" (kept as text), C4="hii"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 was stored as text "30"; the update stores it as a real number.
$ws.Range("B3").Value = 30

# New row 4.
$ws.Range("A4").Value = "hi"

# B4 must stay text "5" (not auto-converted to a number), so force it with
# a leading apostrophe (classic "store as text" trick), then reset the
# cell style back to Normal so no extra quote-prefix formatting lingers.
$ws.Range("B4").Value = "'5"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").Value = "hii"
